$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-8 entirely, shifting cells up so the used range shrinks to A1:E2.
$ws.Rows("3:8").Delete()

# Update the remaining row-2 data (the single winner shown in the popup export).
$ws.Range("A2").Value = 112
# Force B2 to stay text (matches the rest of column B) instead of Excel's
# default numeric auto-detection, then drop the quote-prefix format it
# picks up so the cell keeps the sheet's plain default style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "112"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "Marciana Garay."
$ws.Range("D2").Value = "si"
$ws.Range("E2").Value = "Ganador de Gs. 1.000.000"
